$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "Created a comparison sheet https://github.com/frenziopen/FrenziTech/blob/main/Documentation/Comparison%20Sheet.xlsx`nCreated comparison on NetaFim and wildeye devices and still working further.`nPurchased components and the bill of material is attached https://github.com/frenziopen/FrenziTech/blob/main/Hardware/Components/WhatsApp%20Image%202023-03-09%20at%201.38.14%20PM.jpeg?raw=true"
$ws.Range("C11").WrapText = $true

$ws.Range("C12").Value = "Prepared the hardware of a first kit."

$ws.Rows.Item(11).RowHeight = 114.75
